# issue #5: stock data from json to db
#
# The "股票" (stock) sheet gains two new columns inserted in the middle
# (category, right after property_category / before date) plus two new
# trailing columns (source_file, index) that record where each row of
# stock data originally came from when it was loaded from JSON into the
# database-backed pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column before the existing "date" column (column I)
# so the layout becomes:
#   ... H=property_category, I=<new, category>, J=date, K=legislator_name,
#   L=legislator_id
$ws.Columns.Item(9).Insert()

# Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data rows: every stock row in this export came from the same normalized
# "normal" category, the same source JSON temp file, and keeps its
# original row index (mirrored from column A) for traceability.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpf6f41"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
